$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 84 (shifts existing rows 84..145 down to 85..146,
# carrying their formatting along, matching the target dimension A1:R146).
$ws.Rows("84:84").Insert()

# Populate the freshly-inserted row 84 with the new weekly price record
# (same market/category template as its neighbours, new date + volume + price).
$ws.Cells.Item(84, 1).Value = 5
$ws.Cells.Item(84, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(84, 3).Value = "Maule"
$ws.Cells.Item(84, 4).Value = 44603
$ws.Cells.Item(84, 5).Value = 7
$ws.Cells.Item(84, 6).Value = 100112031
$ws.Cells.Item(84, 7).Value = "Poroto verde"
$ws.Cells.Item(84, 8).Value = "Sin especificar"
$ws.Cells.Item(84, 9).Value = "Primera"
$ws.Cells.Item(84, 10).Value = 100
$ws.Cells.Item(84, 11).Value = 30000
$ws.Cells.Item(84, 12).Value = 30000
$ws.Cells.Item(84, 13).Value = 30000
$ws.Cells.Item(84, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(84, 15).Value = "Región del Maule"
$ws.Cells.Item(84, 16).Value = 1200
$ws.Cells.Item(84, 17).Value = 25
$ws.Cells.Item(84, 18).Value = "Hortaliza"
